# New test added for login also done some config changes
#
# Creds sheet: replace the single hyperlinked "standard_user" credential
# column with a two-column Data Field / Data table (UserName/Password),
# remove the mailto hyperlink, and bold the header row.
#
# New sheet: keep its two values (Arise / 12345) but refresh the view
# (selection + active tab moves to Creds).

$wb = $excel.ActiveWorkbook

$credsWs = $wb.Worksheets.Item("Creds")
$newWs   = $wb.Worksheets.Item("New")

# --- Creds sheet: drop the old hyperlinked single-column layout ---
[void]$credsWs.Hyperlinks.Delete()
[void]$credsWs.Cells.Clear()

# Rebuild as a Data Field / Data table.
$credsWs.Range("B1").Value = "Data"
$credsWs.Range("A1").Value = "Data Field"
$credsWs.Range("A2").Value = "UserName"
$credsWs.Range("B2").Value = "standard_user"
$credsWs.Range("A3").Value = "Password"
$credsWs.Range("B3").Value = "secret_sauce"

# Bold header row.
$credsWs.Range("A1:B1").Font.Bold = $true

# Size the new Data column to fit its contents.
$credsWs.Columns.Item(2).ColumnWidth = 11.6

# --- New sheet: keep the existing values, just refresh the selection ---
$newWs.Range("A6").Select()

# --- Creds becomes the active/selected sheet & cell ---
$credsWs.Range("A3").Select()
$credsWs.Activate()

Write-Host "edit applied"
